$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '24.756.03'
$ws.Range("E2").Value = "'" + '  -0.68%  '
$ws.Range("D3").Value = "'" + '1.679.44'
$ws.Range("E3").Value = "'" + '  -1.81%  '
$ws.Range("D4").Value = "'" + '1.000'
$ws.Range("E4").Value = "'" + '  -0.10%  '
$ws.Range("D5").Value = "'" + '315.07'
$ws.Range("E5").Value = "'" + '  -0.79%  '
$ws.Range("E6").Value = "'" + '  +0.17%  '
$ws.Range("D7").Value = "'" + '0.3923'
$ws.Range("E7").Value = "'" + '  -3.09%  '
$ws.Range("D8").Value = "'" + '0.3961'
$ws.Range("E8").Value = "'" + '  -3.05%  '
$ws.Range("B9").Value = "'" + 'OKB'
$ws.Range("C9").Value = "'" + 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = "'" + '52.12'
$ws.Range("E9").Value = "'" + '  -3.58%  '
$ws.Range("B10").Value = "'" + 'BinanceUSD'
$ws.Range("C10").Value = "'" + 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").Value = "'" + '1.000'
$ws.Range("E10").Value = "'" + '  -0.11%  '
$ws.Range("D11").Value = "'" + '1.402'
$ws.Range("E11").Value = "'" + '  -5.44%  '
$ws.Range("D12").Value = "'" + '0.08647'
$ws.Range("E12").Value = "'" + '  -2.33%  '
$ws.Range("D13").Value = "'" + '25.28'
$ws.Range("E13").Value = "'" + '  -4.50%  '
$ws.Range("D14").Value = "'" + '7.328'
$ws.Range("E14").Value = "'" + '  -2.36%  '
$ws.Range("D15").Value = "'" + '0.00001318'
$ws.Range("E15").Value = "'" + '  -3.15%  '
$ws.Range("D16").Value = "'" + '7.762'
$ws.Range("E16").Value = "'" + '  -4.93%  '
$ws.Range("D17").Value = "'" + '1.662.20'
$ws.Range("E17").Value = "'" + '  -3.13%  '
$ws.Range("D18").Value = "'" + '93.80'
$ws.Range("E18").Value = "'" + '  -3.62%  '
$ws.Range("D19").Value = "'" + '0.07074'
$ws.Range("E19").Value = "'" + '  -1.31%  '
$ws.Range("D20").Value = "'" + '20.46'
$ws.Range("E20").Value = "'" + '  -3.65%  '
$ws.Range("D21").Value = "'" + '7.077'
$ws.Range("E21").Value = "'" + '  -2.91%  '
$ws.Range("D22").Value = "'" + '1.003'
$ws.Range("E22").Value = "'" + '  +0.25%  '
$ws.Range("D23").Value = "'" + '13.96'
$ws.Range("E23").Value = "'" + '  -3.21%  '
$ws.Range("D24").Value = "'" + '24.748.74'
$ws.Range("E24").Value = "'" + '  -0.72%  '
$ws.Range("D25").Value = "'" + '2.348'
$ws.Range("E25").Value = "'" + '  +0.86%  '
$ws.Range("D26").Value = "'" + '2.785'
$ws.Range("E26").Value = "'" + '  -4.99%  '
$ws.Range("D27").Value = "'" + '23.37'
$ws.Range("E27").Value = "'" + '  -0.15%  '
$ws.Range("D28").Value = "'" + '162.07'
$ws.Range("E28").Value = "'" + '  -2.98%  '
$ws.Range("D29").Value = "'" + '5.764'
$ws.Range("E29").Value = "'" + '  -7.56%  '
$ws.Range("D30").Value = "'" + '147.17'
$ws.Range("E30").Value = "'" + '  +0.62%  '
$ws.Range("D31").Value = "'" + '7.873'
$ws.Range("E31").Value = "'" + '  -6.98%  '
$ws.Range("D32").Value = "'" + '2.446'
$ws.Range("E32").Value = "'" + '  +8.77%  '
$ws.Range("D33").Value = "'" + '1.832.03'
$ws.Range("E33").Value = "'" + '  -3.47%  '
$ws.Range("D34").Value = "'" + '0.08436'
$ws.Range("E34").Value = "'" + '  -4.62%  '
$ws.Range("D35").Value = "'" + '0.03041'
$ws.Range("E35").Value = "'" + '  -4.91%  '
$ws.Range("D36").Value = "'" + '6.939'
$ws.Range("E36").Value = "'" + '  -4.31%  '
$ws.Range("D37").Value = "'" + '0.2822'
$ws.Range("E37").Value = "'" + '  -2.01%  '
$ws.Range("D38").Value = "'" + '0.9983'
$ws.Range("E38").Value = "'" + '  -3.17%  '
$ws.Range("D39").Value = "'" + '0.09505'
$ws.Range("E39").Value = "'" + '  +1.48%  '
$ws.Range("D40").Value = "'" + '10.63'
$ws.Range("E40").Value = "'" + '  -2.53%  '
$ws.Range("D41").Value = "'" + '1.526'
$ws.Range("E41").Value = "'" + '  +3.69%  '
$ws.Range("D42").Value = "'" + '0.7942'
$ws.Range("E42").Value = "'" + '  -6.24%  '
$ws.Range("D43").Value = "'" + '13.57'
$ws.Range("E43").Value = "'" + '  -4.43%  '
$ws.Range("D44").Value = "'" + '16.58'
$ws.Range("E44").Value = "'" + '  -5.32%  '
$ws.Range("D45").Value = "'" + '0.7149'
$ws.Range("E45").Value = "'" + '  -3.97%  '
$ws.Range("D46").Value = "'" + '2.569'
$ws.Range("E46").Value = "'" + '  -5.60%  '
$ws.Range("D47").Value = "'" + '4.197'
$ws.Range("E47").Value = "'" + '  -1.05%  '
$ws.Range("D48").Value = "'" + '0.08677'
$ws.Range("E48").Value = "'" + '  +3.25%  '
$ws.Range("D49").Value = "'" + '1.001'
$ws.Range("E49").Value = "'" + '  +0.22%  '
$ws.Range("E50").Value = "'" + '  -4.47%  '
$ws.Range("D51").Value = "'" + '138.03'
$ws.Range("E51").Value = "'" + '  -2.73%  '
